$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force Excel to store the value as text (not auto-convert to a number),
    # while keeping the cell's style/format identical to before (no style index).
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

# Row 2
Set-TextValue "D2" "68.777.09"
$ws.Range("E2").Value = "  -0.48%  "

# Row 3
Set-TextValue "D3" "3.862.64"
$ws.Range("E3").Value = "  +3.04%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
Set-TextValue "D5" "600.27"
$ws.Range("E5").Value = "  -0.19%  "

# Row 6
Set-TextValue "D6" "162.14"
$ws.Range("E6").Value = "  -2.89%  "

# Row 7
Set-TextValue "D7" "3.861.74"
$ws.Range("E7").Value = "  +3.02%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("E9").Value = "  -1.82%  "

# Row 10
$ws.Range("E10").Value = "  -0.97%  "

# Row 11
$ws.Range("E11").Value = "  -1.15%  "

# Row 12
$ws.Range("E12").Value = "  -0.38%  "

# Row 13
Set-TextValue "D13" "36.85"
$ws.Range("E13").Value = "  -2.92%  "

# Row 14
$ws.Range("E14").Value = "  -2.12%  "

# Row 15
Set-TextValue "D15" "4.512.99"
$ws.Range("E15").Value = "  +3.17%  "

# Row 16
Set-TextValue "D16" "3.804.03"
$ws.Range("E16").Value = "  +1.46%  "

# Row 17
Set-TextValue "D17" "68.961.42"
$ws.Range("E17").Value = "  -0.20%  "

# Row 18
$ws.Range("E18").Value = "  +2.66%  "

# Row 19
$ws.Range("E19").Value = "  -0.41%  "

# Row 20
Set-TextValue "D20" "11.38"
$ws.Range("E20").Value = "  +3.00%  "

# Row 21
Set-TextValue "D21" "17.14"
$ws.Range("E21").Value = "  -1.52%  "

# Row 22
Set-TextValue "D22" "483.37"
$ws.Range("E22").Value = "  -1.86%  "

# Row 23
Set-TextValue "D23" "0.717"
$ws.Range("E23").Value = "  -1.39%  "

# Row 24
$ws.Range("E24").Value = "  +6.51%  "

# Row 25
Set-TextValue "D25" "83.91"
$ws.Range("E25").Value = "  -1.12%  "

# Row 26
$ws.Range("E26").Value = "  -2.92%  "

# Row 27
Set-TextValue "D27" "12.08"
$ws.Range("E27").Value = "  -1.55%  "

# Row 28
$ws.Range("E28").Value = "  -0.07%  "

# Row 29
Set-TextValue "D29" "9.94"
$ws.Range("E29").Value = "  -1.08%  "

# Row 31
Set-TextValue "D31" "4.015.88"
$ws.Range("E31").Value = "  +3.10%  "

# Row 32
Set-TextValue "D32" "7.86"
$ws.Range("E32").Value = "  -3.16%  "

# Row 33
Set-TextValue "D33" "32.20"

# Row 34
$ws.Range("E34").Value = "  -4.31%  "

# Row 35
Set-TextValue "D35" "3.810.44"

# Row 36
$ws.Range("E36").Value = "  -1.70%  "

# Row 37
$ws.Range("E37").Value = "  +1.98%  "

# Row 38
$ws.Range("E38").Value = "  +1.50%  "

# Row 39
Set-TextValue "D39" "5.86"
$ws.Range("E39").Value = "  -1.60%  "

# Row 40
$ws.Range("E40").Value = "  +0.03%  "

# Row 41
Set-TextValue "D41" "0.317"
$ws.Range("E41").Value = "  -2.70%  "

# Row 42
$ws.Range("E42").Value = "  -2.10%  "

# Row 43
$ws.Range("E43").Value = "  +1.55%  "

# Row 44
$ws.Range("E44").Value = "  -0.31%  "

# Row 45
$ws.Range("E45").Value = "  -0.08%  "

# Row 47
Set-TextValue "D47" "8.37"
$ws.Range("E47").Value = "  -1.06%  "

# Row 48
Set-TextValue "D48" "143.65"
$ws.Range("E48").Value = "  +1.66%  "

# Row 49
Set-TextValue "D49" "2.838.48"
$ws.Range("E49").Value = "  +1.74%  "

# Row 50
$ws.Range("E50").Value = "  +1.09%  "

# Row 51
Set-TextValue "D51" "25.83"
$ws.Range("E51").Value = "  +12.67%  "
